# Update odds figures for Jogos_da_Semana_FlashScore_2025-03-21.xlsx (Sheet1)
# Applies the updated odd values for rows 2, 3, 5, 6, 7, 11, 12, 13, 15, 16
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 2.12   # G2
$ws.Cells.Item(2, 8).Value = 2.75   # H2
$ws.Cells.Item(2, 9).Value = 3.95   # I2
$ws.Cells.Item(2, 10).Value = 2.85   # J2
$ws.Cells.Item(2, 11).Value = 1.83   # K2
$ws.Cells.Item(2, 12).Value = 4.5   # L2
$ws.Cells.Item(2, 13).Value = 1.39   # M2
$ws.Cells.Item(2, 14).Value = 2.57   # N2
$ws.Cells.Item(2, 15).Value = 2.12   # O2
$ws.Cells.Item(2, 16).Value = 1.57   # P2
$ws.Cells.Item(2, 17).Value = 3.5   # Q2
$ws.Cells.Item(2, 18).Value = 1.21   # R2
$ws.Cells.Item(2, 19).Value = 1.53   # S2
$ws.Cells.Item(2, 20).Value = 2.2   # T2
$ws.Cells.Item(2, 21).Value = 1.78   # U2
$ws.Cells.Item(2, 22).Value = 1.82   # V2
$ws.Cells.Item(2, 23).Value = 6.4   # W2
$ws.Cells.Item(2, 24).Value = 9.75   # X2
$ws.Cells.Item(2, 25).Value = 8.5   # Y2
$ws.Cells.Item(2, 26).Value = 22   # Z2
$ws.Cells.Item(2, 27).Value = 18.5   # AA2
$ws.Cells.Item(2, 28).Value = 30   # AB2
$ws.Cells.Item(2, 29).Value = 6.9   # AC2
$ws.Cells.Item(2, 30).Value = 5.4   # AD2
$ws.Cells.Item(2, 31).Value = 13.5   # AE2
$ws.Cells.Item(2, 32).Value = 70   # AF2
$ws.Cells.Item(2, 33).Value = 600   # AG2
$ws.Cells.Item(2, 34).Value = 9.75   # AH2
$ws.Cells.Item(2, 35).Value = 22   # AI2
$ws.Cells.Item(2, 36).Value = 13   # AJ2
$ws.Cells.Item(2, 38).Value = 40   # AL2
$ws.Cells.Item(2, 39).Value = 45   # AM2
# Row 3
$ws.Cells.Item(3, 8).Value = 3.6   # H3
$ws.Cells.Item(3, 9).Value = 4.75   # I3
$ws.Cells.Item(3, 10).Value = 2.18   # J3
$ws.Cells.Item(3, 14).Value = 2.87   # N3
$ws.Cells.Item(3, 15).Value = 1.93   # O3
$ws.Cells.Item(3, 16).Value = 1.7   # P3
$ws.Cells.Item(3, 17).Value = 3.15   # Q3
$ws.Cells.Item(3, 18).Value = 1.26   # R3
$ws.Cells.Item(3, 22).Value = 1.7   # V3
$ws.Cells.Item(3, 23).Value = 6.1   # W3
$ws.Cells.Item(3, 29).Value = 9.25   # AC3
$ws.Cells.Item(3, 32).Value = 100   # AF3
$ws.Cells.Item(3, 33).Value = 900   # AG3
$ws.Cells.Item(3, 34).Value = 11.75   # AH3
$ws.Cells.Item(3, 39).Value = 55   # AM3
# Row 5
$ws.Cells.Item(5, 7).Value = 2.5   # G5
$ws.Cells.Item(5, 8).Value = 2.85   # H5
$ws.Cells.Item(5, 9).Value = 2.92   # I5
$ws.Cells.Item(5, 10).Value = 3.25   # J5
$ws.Cells.Item(5, 11).Value = 1.87   # K5
$ws.Cells.Item(5, 12).Value = 3.6   # L5
$ws.Cells.Item(5, 13).Value = 1.52   # M5
$ws.Cells.Item(5, 14).Value = 2.22   # N5
$ws.Cells.Item(5, 15).Value = 2.47   # O5
$ws.Cells.Item(5, 16).Value = 1.42   # P5
$ws.Cells.Item(5, 17).Value = 4.3   # Q5
$ws.Cells.Item(5, 18).Value = 1.14   # R5
$ws.Cells.Item(5, 19).Value = 1.55   # S5
$ws.Cells.Item(5, 20).Value = 2.15   # T5
$ws.Cells.Item(5, 21).Value = 2.07   # U5
$ws.Cells.Item(5, 22).Value = 1.6   # V5
$ws.Cells.Item(5, 23).Value = 5.9   # W5
$ws.Cells.Item(5, 24).Value = 10.75   # X5
$ws.Cells.Item(5, 25).Value = 10.5   # Y5
$ws.Cells.Item(5, 26).Value = 28   # Z5
$ws.Cells.Item(5, 27).Value = 28   # AA5
$ws.Cells.Item(5, 28).Value = 50   # AB5
$ws.Cells.Item(5, 29).Value = 6.1   # AC5
$ws.Cells.Item(5, 30).Value = 5.8   # AD5
$ws.Cells.Item(5, 31).Value = 18.5   # AE5
$ws.Cells.Item(5, 32).Value = 120   # AF5
$ws.Cells.Item(5, 34).Value = 6.8   # AH5
$ws.Cells.Item(5, 35).Value = 13.5   # AI5
$ws.Cells.Item(5, 36).Value = 11.25   # AJ5
$ws.Cells.Item(5, 37).Value = 37   # AK5
$ws.Cells.Item(5, 38).Value = 32   # AL5
$ws.Cells.Item(5, 39).Value = 50   # AM5
$ws.Cells.Item(5, 40).Value = 1.1   # AN5
$ws.Cells.Item(5, 41).Value = 6.2   # AO5
# Row 6
$ws.Cells.Item(6, 7).Value = 2.32   # G6
$ws.Cells.Item(6, 9).Value = 2.77   # I6
$ws.Cells.Item(6, 10).Value = 2.87   # J6
$ws.Cells.Item(6, 11).Value = 2.15   # K6
$ws.Cells.Item(6, 12).Value = 3.3   # L6
$ws.Cells.Item(6, 13).Value = 1.32   # M6
$ws.Cells.Item(6, 14).Value = 2.82   # N6
$ws.Cells.Item(6, 15).Value = 1.93   # O6
$ws.Cells.Item(6, 16).Value = 1.7   # P6
$ws.Cells.Item(6, 22).Value = 1.83   # V6
$ws.Cells.Item(6, 26).Value = 23   # Z6
$ws.Cells.Item(6, 30).Value = 6.5   # AD6
$ws.Cells.Item(6, 31).Value = 15   # AE6
$ws.Cells.Item(6, 33).Value = 600   # AG6
$ws.Cells.Item(6, 35).Value = 13.5   # AI6
$ws.Cells.Item(6, 36).Value = 10.5   # AJ6
$ws.Cells.Item(6, 37).Value = 32   # AK6
# Row 7
$ws.Cells.Item(7, 7).Value = 1.82   # G7
$ws.Cells.Item(7, 9).Value = 3.8   # I7
$ws.Cells.Item(7, 10).Value = 2.4   # J7
$ws.Cells.Item(7, 12).Value = 4.15   # L7
$ws.Cells.Item(7, 17).Value = 2.8   # Q7
$ws.Cells.Item(7, 18).Value = 1.33   # R7
$ws.Cells.Item(7, 22).Value = 1.88   # V7
$ws.Cells.Item(7, 23).Value = 7.4   # W7
$ws.Cells.Item(7, 25).Value = 8.25   # Y7
$ws.Cells.Item(7, 26).Value = 14.5   # Z7
$ws.Cells.Item(7, 27).Value = 14   # AA7
$ws.Cells.Item(7, 34).Value = 11.5   # AH7
$ws.Cells.Item(7, 35).Value = 21   # AI7
$ws.Cells.Item(7, 36).Value = 13   # AJ7
$ws.Cells.Item(7, 37).Value = 55   # AK7
$ws.Cells.Item(7, 38).Value = 35   # AL7
# Row 11
$ws.Cells.Item(11, 7).Value = 2.3   # G11
$ws.Cells.Item(11, 8).Value = 2.8   # H11
$ws.Cells.Item(11, 9).Value = 3.4   # I11
$ws.Cells.Item(11, 10).Value = 3.1   # J11
$ws.Cells.Item(11, 11).Value = 1.83   # K11
$ws.Cells.Item(11, 13).Value = 1.58   # M11
$ws.Cells.Item(11, 14).Value = 2.2   # N11
$ws.Cells.Item(11, 15).Value = 2.75   # O11
$ws.Cells.Item(11, 16).Value = 1.4   # P11
$ws.Cells.Item(11, 17).Value = 5.5   # Q11
$ws.Cells.Item(11, 18).Value = 1.1   # R11
$ws.Cells.Item(11, 19).Value = 1.62   # S11
$ws.Cells.Item(11, 20).Value = 2.2   # T11
$ws.Cells.Item(11, 21).Value = 2.25   # U11
$ws.Cells.Item(11, 22).Value = 1.57   # V11
$ws.Cells.Item(11, 23).Value = 6   # W11
$ws.Cells.Item(11, 25).Value = 11   # Y11
$ws.Cells.Item(11, 27).Value = 23   # AA11
$ws.Cells.Item(11, 29).Value = 5.5   # AC11
$ws.Cells.Item(11, 31).Value = 21   # AE11
$ws.Cells.Item(11, 32).Value = 81   # AF11
$ws.Cells.Item(11, 34).Value = 7.5   # AH11
$ws.Cells.Item(11, 36).Value = 15   # AJ11
$ws.Cells.Item(11, 39).Value = 51   # AM11
$ws.Cells.Item(11, 40).Value = 1.11   # AN11
$ws.Cells.Item(11, 41).Value = 5   # AO11
# Row 12
$ws.Cells.Item(12, 7).Value = 3.4   # G12
$ws.Cells.Item(12, 8).Value = 2.9   # H12
$ws.Cells.Item(12, 9).Value = 2.3   # I12
$ws.Cells.Item(12, 10).Value = 4   # J12
$ws.Cells.Item(12, 11).Value = 1.91   # K12
$ws.Cells.Item(12, 12).Value = 3.2   # L12
$ws.Cells.Item(12, 13).Value = 1.47   # M12
$ws.Cells.Item(12, 14).Value = 2.5   # N12
$ws.Cells.Item(12, 15).Value = 2.6   # O12
$ws.Cells.Item(12, 16).Value = 1.48   # P12
$ws.Cells.Item(12, 17).Value = 5   # Q12
$ws.Cells.Item(12, 18).Value = 1.13   # R12
$ws.Cells.Item(12, 19).Value = 1.57   # S12
$ws.Cells.Item(12, 20).Value = 2.25   # T12
$ws.Cells.Item(12, 23).Value = 8   # W12
$ws.Cells.Item(12, 24).Value = 15   # X12
$ws.Cells.Item(12, 26).Value = 34   # Z12
$ws.Cells.Item(12, 29).Value = 6.5   # AC12
$ws.Cells.Item(12, 35).Value = 10   # AI12
$ws.Cells.Item(12, 36).Value = 10   # AJ12
$ws.Cells.Item(12, 37).Value = 21   # AK12
$ws.Cells.Item(12, 38).Value = 23   # AL12
$ws.Cells.Item(12, 40).Value = 1.08   # AN12
$ws.Cells.Item(12, 41).Value = 6.5   # AO12
$ws.Cells.Item(12, 42).Value = 1.98   # AP12
$ws.Cells.Item(12, 43).Value = 1.88   # AQ12
# Row 13
$ws.Cells.Item(13, 7).Value = 3.7   # G13
$ws.Cells.Item(13, 8).Value = 3.55   # H13
$ws.Cells.Item(13, 9).Value = 1.9   # I13
$ws.Cells.Item(13, 10).Value = 4.1   # J13
$ws.Cells.Item(13, 11).Value = 2.18   # K13
$ws.Cells.Item(13, 12).Value = 2.45   # L13
$ws.Cells.Item(13, 13).Value = 1.24   # M13
$ws.Cells.Item(13, 14).Value = 3.65   # N13
$ws.Cells.Item(13, 15).Value = 1.72   # O13
$ws.Cells.Item(13, 16).Value = 2   # P13
$ws.Cells.Item(13, 17).Value = 2.72   # Q13
$ws.Cells.Item(13, 18).Value = 1.4   # R13
$ws.Cells.Item(13, 19).Value = 1.36   # S13
$ws.Cells.Item(13, 20).Value = 2.9   # T13
$ws.Cells.Item(13, 23).Value = 12   # W13
$ws.Cells.Item(13, 24).Value = 21   # X13
$ws.Cells.Item(13, 25).Value = 12.5   # Y13
$ws.Cells.Item(13, 26).Value = 55   # Z13
$ws.Cells.Item(13, 27).Value = 32   # AA13
$ws.Cells.Item(13, 28).Value = 35   # AB13
$ws.Cells.Item(13, 29).Value = 7.9   # AC13
$ws.Cells.Item(13, 30).Value = 6.9   # AD13
$ws.Cells.Item(13, 31).Value = 13   # AE13
$ws.Cells.Item(13, 34).Value = 8.5   # AH13
$ws.Cells.Item(13, 35).Value = 10   # AI13
$ws.Cells.Item(13, 36).Value = 8.25   # AJ13
$ws.Cells.Item(13, 37).Value = 17   # AK13
$ws.Cells.Item(13, 38).Value = 14   # AL13
$ws.Cells.Item(13, 40).Value = 1.05   # AN13
$ws.Cells.Item(13, 41).Value = 7.9   # AO13
# Row 15
$ws.Cells.Item(15, 11).Value = 2.2   # K15
$ws.Cells.Item(15, 12).Value = 3.05   # L15
$ws.Cells.Item(15, 14).Value = 3.4   # N15
$ws.Cells.Item(15, 16).Value = 1.9   # P15
$ws.Cells.Item(15, 18).Value = 1.35   # R15
$ws.Cells.Item(15, 20).Value = 2.87   # T15
$ws.Cells.Item(15, 32).Value = 60   # AF15
# Row 16
$ws.Cells.Item(16, 7).Value = 2.32   # G16
$ws.Cells.Item(16, 8).Value = 3   # H16
$ws.Cells.Item(16, 10).Value = 2.95   # J16
$ws.Cells.Item(16, 11).Value = 1.98   # K16
$ws.Cells.Item(16, 12).Value = 3.75   # L16
$ws.Cells.Item(16, 13).Value = 1.42   # M16
$ws.Cells.Item(16, 14).Value = 2.65   # N16
$ws.Cells.Item(16, 15).Value = 2.25   # O16
$ws.Cells.Item(16, 16).Value = 1.57   # P16
$ws.Cells.Item(16, 17).Value = 3.9   # Q16
$ws.Cells.Item(16, 18).Value = 1.21   # R16
$ws.Cells.Item(16, 19).Value = 1.47   # S16
$ws.Cells.Item(16, 20).Value = 2.5   # T16
$ws.Cells.Item(16, 21).Value = 1.9   # U16
$ws.Cells.Item(16, 22).Value = 1.8   # V16
$ws.Cells.Item(16, 23).Value = 6.6   # W16
$ws.Cells.Item(16, 24).Value = 10.5   # X16
$ws.Cells.Item(16, 25).Value = 9.25   # Y16
$ws.Cells.Item(16, 27).Value = 21   # AA16
$ws.Cells.Item(16, 28).Value = 35   # AB16
$ws.Cells.Item(16, 29).Value = 6   # AC16
$ws.Cells.Item(16, 31).Value = 15   # AE16
$ws.Cells.Item(16, 32).Value = 80   # AF16
$ws.Cells.Item(16, 33).Value = 800   # AG16
$ws.Cells.Item(16, 34).Value = 8.25   # AH16
$ws.Cells.Item(16, 35).Value = 15.5   # AI16
$ws.Cells.Item(16, 36).Value = 11.25   # AJ16
$ws.Cells.Item(16, 38).Value = 30   # AL16
$ws.Cells.Item(16, 39).Value = 40   # AM16
$ws.Cells.Item(16, 40).Value = 1.1   # AN16
$ws.Cells.Item(16, 41).Value = 6   # AO16
